$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Add new row 19 data (mail log entry)
$ws.Cells.Item(19, 1).Value = "Kun je 5 rollen afplaktape bestellen?"
$ws.Cells.Item(19, 2).Value = "MailMind Test <mailmind.test@zohomail.eu>"
$ws.Cells.Item(19, 3).Value = "Hoi Johan,`n Zou je 5 rollen afplaktape kunnen bestellen voor de schildersafdeling?`nDank je wel!`n — Marco`nSent using {0}"
$ws.Cells.Item(19, 4).Value = "Bestelling / Levering"
$ws.Cells.Item(19, 5).Value = "Beste Marco,`nBedankt voor je e-mail. We zullen direct 5 rollen afplaktape bestellen voor de schildersafdeling. Zodra de bestelling is geplaatst, ontvang je hier een bevestiging van.`nMet vriendelijke groet,`nJohan"
$ws.Cells.Item(19, 6).Value = "2025-06-26 21:25:07"
$ws.Cells.Item(19, 7).Value = "Ja"
$ws.Cells.Item(19, 8).Value = "Nee"
$ws.Cells.Item(19, 9).Value = "Ja"

# Reset row height to default so multi-line content doesn't leave a stray custom row height
$ws.Rows.Item(19).EntireRow.AutoFit()

# Extend conditional formatting ranges so they also cover the newly added row 19
$ws.Range("D2:D18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D19"))
$ws.Range("G2:G18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G19"))
$ws.Range("H2:H18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H19"))
$ws.Range("I2:I18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I19"))

# Update Dashboard summary count for "Bestelling / Levering" (13 -> 14)
$wsDash = $wb.Worksheets.Item("Dashboard")
$wsDash.Cells.Item(2, 2).Value = 14
